$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded for Cilantro at "Terminal La
# Palmera de La Serena". It belongs right after the existing row 17, so
# insert a fresh row at position 18 (this pushes the former rows 18-79 down
# to 19-80, carrying their values/styles with them - including the row that
# lands at 80, which is brand new in the sheet).
$ws.Rows.Item(18).Insert()

# Fill in the new row 18 with the new observation. All the "descriptive"
# columns (market, region, product, etc.) are constant for every row in this
# sheet, so reuse them; only the date (D), volume (J) and the
# min/max/avg/kg-price columns (K/L/M/P) actually vary per row - and here
# they match what used to be in row 18 before the shift (i.e. same price
# bucket, new date/volume).
$ws.Cells.Item(18, 1).Value = 8
$ws.Cells.Item(18, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(18, 3).Value = "Coquimbo"
$ws.Cells.Item(18, 4).Value = 44453
$ws.Cells.Item(18, 5).Value = 4
$ws.Cells.Item(18, 6).Value = 100112040
$ws.Cells.Item(18, 7).Value = "Cilantro"
$ws.Cells.Item(18, 8).Value = "Sin especificar"
$ws.Cells.Item(18, 9).Value = "Primera"
$ws.Cells.Item(18, 10).Value = 3300
$ws.Cells.Item(18, 11).Value = 2000
$ws.Cells.Item(18, 12).Value = 2500
$ws.Cells.Item(18, 13).Value = 2250
$ws.Cells.Item(18, 14).Value = "`$/atado 1 a 1,5 kilos"
$ws.Cells.Item(18, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(18, 16).Value = 1500
$ws.Cells.Item(18, 17).Value = 1.5
$ws.Cells.Item(18, 18).Value = "Hortaliza"
